$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from inside "anstranslation" (mid-word) to a
#    zero-width bookmark right after "Run the service" (end of the
#    "Step 4: Run the service" heading paragraph).
# ---------------------------------------------------------------------------

# Locate the heading paragraph so we get its live end offset (robust to any
# earlier edits in the document).
$headingEnd = -1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Step 4: Run the service") {
        $headingEnd = $p.Range.End
    }
}

if ($headingEnd -ge 0) {
    # The engine's Bookmarks.Add mis-resolves a *collapsed* range sitting
    # exactly one character before a paragraph end (i.e. start == end ==
    # paragraph.End - 1): it wrongly anchors at the very start of the
    # document instead of the requested offset. Work around this by
    # temporarily inserting a one-character sentinel right after "service"
    # so the insertion point we want is no longer the paragraph's last
    # character, add the bookmark there (now safely mid-paragraph), then
    # delete the sentinel again.
    $insertionPoint = $headingEnd - 1
    $sentinelRange = $d.Range($insertionPoint, $insertionPoint)
    $sentinelRange.InsertAfter("X")

    $bmRange = $d.Range($insertionPoint, $insertionPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $sentinelDelRange = $d.Range($insertionPoint, $insertionPoint + 1)
    $sentinelDelRange.Delete()
}

# ---------------------------------------------------------------------------
# 2) Update the service IP address and merge the run containing
#    "Service will be available at: http://" + old IP + " and http://" into
#    a single run with the new IP baked in.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Service will be available at: http://34.163.137.242 and http://",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Service will be available at: http://34.163.79.207 and http://", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Merge the old "anstransl" + bookmark + "ation" split runs into a single
#    "anstranslation" run (the stray bookmark that used to live here was
#    already relocated in step 1, so this just re-merges the two runs).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "anstranslation", $false, $false, $false, $false, $false, $true, 1,
    $false, "anstranslation", 2) | Out-Null

Write-Output "done"
